$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product1")

# Replace the volatile formula in D7 with a static text value
$ws.Range("D7").Value = "sample.pdf"

# Update the active selection to D8
$ws.Range("D8").Select()
